# Updates cryptos price/volume table cells to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellRefs = @(
    'D2',
    'E2',
    'D3',
    'E3',
    'D4',
    'E4',
    'D5',
    'E5',
    'D6',
    'E6',
    'D7',
    'E7',
    'E8',
    'D9',
    'E9',
    'D10',
    'E10',
    'D11',
    'E11',
    'D12',
    'E12',
    'D13',
    'E13',
    'D14',
    'E14',
    'D15',
    'E15',
    'E16',
    'D17',
    'E17',
    'D18',
    'E18',
    'E19',
    'D20',
    'E20',
    'D21',
    'E21',
    'D22',
    'E22',
    'D23',
    'E23',
    'B24',
    'C24',
    'D24',
    'E24',
    'B25',
    'C25',
    'D25',
    'E25',
    'D26',
    'E26',
    'D27',
    'E27',
    'D28',
    'D29',
    'E29',
    'D30',
    'E30',
    'D31',
    'E31',
    'E32',
    'D33',
    'E33',
    'D34',
    'E34',
    'D35',
    'E35',
    'D36',
    'E36',
    'D37',
    'E37',
    'D38',
    'E38',
    'D39',
    'E39',
    'D40',
    'E40',
    'D41',
    'E41',
    'D42',
    'E42',
    'D43',
    'E43',
    'D44',
    'E44',
    'D45',
    'E45',
    'D46',
    'E46',
    'D47',
    'E47',
    'D48',
    'E48',
    'D49',
    'E49',
    'D50',
    'E50',
    'D51',
    'E51'
)

$newValues = @(
    '27.152.30',
    '  -0.15%  ',
    '1.898.61',
    '  -0.28%  ',
    '1.004',
    '  +0.39%  ',
    '306.95',
    '  +0.19%  ',
    '1.003',
    '  +0.30%  ',
    '0.5235',
    '  -0.29%  ',
    '  +0.84%  ',
    '0.07283',
    '  +0.38%  ',
    '21.37',
    '  +1.13%  ',
    '0.9029',
    '  +0.41%  ',
    '0.08161',
    '  -3.01%  ',
    '95.36',
    '  +0.71%  ',
    '5.343',
    '  +1.46%  ',
    '1.806.34',
    '  -5.09%  ',
    '  +0.35%  ',
    '0.000008651',
    '  +0.53%  ',
    '14.68',
    '  +0.88%  ',
    '  +0.34%  ',
    '27.188.27',
    '  -0.16%  ',
    '5.096',
    '  +0.73%  ',
    '10.78',
    '  +1.78%  ',
    '6.445',
    '  +0.13%  ',
    'Monero',
    'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr',
    '149.50',
    '  +1.79%  ',
    'LidoDAOToken',
    'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo',
    '2.325',
    '  +2.12%  ',
    '18.20',
    '  +0.21%  ',
    '1.743',
    '  -0.50%  ',
    '115.83',
    '4.824',
    '  +0.41%  ',
    '4.877',
    '  -1.03%  ',
    '0.09219',
    '  -0.70%  ',
    '  -0.32%  ',
    '0.7931',
    '  -2.06%  ',
    '1.222',
    '  -1.17%  ',
    '2.966',
    '  +0.43%  ',
    '3.368',
    '  +0.44%  ',
    '2.655',
    '  +1.72%  ',
    '0.5696',
    '  -0.15%  ',
    '0.01987',
    '  -0.18%  ',
    '1.080',
    '  +0.77%  ',
    '9.005',
    '  +0.42%  ',
    '6.579',
    '  -1.14%  ',
    '116.05',
    '  -1.44%  ',
    '0.1510',
    '  -0.32%  ',
    '0.4882',
    '  +0.81%  ',
    '1.003',
    '  +0.34%  ',
    '10.14',
    '  +0.00%  ',
    '1.625',
    '  +0.83%  ',
    '38.44',
    '  +2.67%  ',
    '63.82',
    '  +0.36%  ',
    '0.05958',
    '  +0.49%  '
)

for ($i = 0; $i -lt $cellRefs.Length; $i++) {
    # Prefix with a literal apostrophe so Excel stores the value as text
    # (avoids values such as "1.004" or "27.152.30" being parsed as numbers),
    # then clear the resulting quote-prefix formatting so no stray style is left behind.
    $ws.Range($cellRefs[$i]).Value = "'" + $newValues[$i]
    $ws.Range($cellRefs[$i]).ClearFormats()
}
